$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving data row's value
$ws.Range("A2").Value = "G00004"

# Remove the now-missing rows (3-8) entirely
$ws.Range("A3:A8").EntireRow.Delete() | Out-Null

# Match the saved selection state
$ws.Range("A3:A8").Select() | Out-Null
